$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.863.94'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '2.115.78'
$ws.Range('D4').NumberFormatLocal = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormatLocal = '@'
$ws.Range('D5').Value = '348.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').NumberFormatLocal = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').NumberFormatLocal = '@'
$ws.Range('D8').Value = '0.4468'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('D9').NumberFormatLocal = '@'
$ws.Range('D9').Value = '54.21'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.49%  '
$ws.Range('D10').NumberFormatLocal = '@'
$ws.Range('D10').Value = '0.09355'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.33%  '
$ws.Range('D11').NumberFormatLocal = '@'
$ws.Range('D11').Value = '1.182'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').NumberFormatLocal = '@'
$ws.Range('D12').Value = '25.22'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('D13').Value = '2.107.42'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormatLocal = '@'
$ws.Range('D14').Value = '8.377'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.99%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormatLocal = '@'
$ws.Range('D15').Value = '6.842'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('E16').Value = '  +3.82%  '
$ws.Range('D17').NumberFormatLocal = '@'
$ws.Range('D17').Value = '0.00001167'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').NumberFormatLocal = '@'
$ws.Range('D19').Value = '21.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.62%  '
$ws.Range('D20').NumberFormatLocal = '@'
$ws.Range('D20').Value = '0.06674'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = '29.916.39'
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('D24').NumberFormatLocal = '@'
$ws.Range('D24').Value = '12.73'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').NumberFormatLocal = '@'
$ws.Range('D25').Value = '2.330'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Value = '2.357.14'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormatLocal = '@'
$ws.Range('D27').Value = '22.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.10%  '
$ws.Range('D28').NumberFormatLocal = '@'
$ws.Range('D28').Value = '2.554'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormatLocal = '@'
$ws.Range('D29').Value = '162.60'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').NumberFormatLocal = '@'
$ws.Range('D30').Value = '134.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormatLocal = '@'
$ws.Range('D31').Value = '1.156'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('D32').NumberFormatLocal = '@'
$ws.Range('D32').Value = '1.793'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.24%  '
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').NumberFormatLocal = '@'
$ws.Range('D34').Value = '6.252'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('D35').NumberFormatLocal = '@'
$ws.Range('D35').Value = '3.972'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormatLocal = '@'
$ws.Range('D36').Value = '6.403'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.62%  '
$ws.Range('D37').NumberFormatLocal = '@'
$ws.Range('D37').Value = '10.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.58%  '
$ws.Range('D38').NumberFormatLocal = '@'
$ws.Range('D38').Value = '0.02596'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.98%  '
$ws.Range('D39').NumberFormatLocal = '@'
$ws.Range('D39').Value = '0.06810'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormatLocal = '@'
$ws.Range('D40').Value = '12.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('D41').NumberFormatLocal = '@'
$ws.Range('D41').Value = '0.7021'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormatLocal = '@'
$ws.Range('D43').Value = '0.2247'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').NumberFormatLocal = '@'
$ws.Range('D44').Value = '0.6852'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.85%  '
$ws.Range('D45').NumberFormatLocal = '@'
$ws.Range('D45').Value = '14.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('D46').NumberFormatLocal = '@'
$ws.Range('D46').Value = '2.363'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.50%  '
$ws.Range('D47').NumberFormatLocal = '@'
$ws.Range('D47').Value = '1.006'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('D48').NumberFormatLocal = '@'
$ws.Range('D48').Value = '3.636'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormatLocal = '@'
$ws.Range('D49').Value = '0.00000000358'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('D50').NumberFormatLocal = '@'
$ws.Range('D50').Value = '1.217'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.21%  '
$ws.Range('E51').Value = '  +0.75%  '
